$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# ------------------------------------------------------------------
# 1) Create new row 17 as a copy of row 3 (values + formatting), BEFORE
#    row 2/3 get their own edits below, so row 17 inherits row 3's
#    original ("Palak"/"abcd1234$"/"IndiaTetherfi"/"Cheker") content.
# ------------------------------------------------------------------
$ws.Range("A3:F3").Copy($ws.Range("A17:F17"))

# ------------------------------------------------------------------
# 2) Row 2 + Row 3: Application URL column (A) now points at the new
#    UI, and loses its live hyperlink - restyle like A16 (hyperlink
#    color, but not an actual hyperlink, no vertical-center alignment).
# ------------------------------------------------------------------
$ws.Range("A2").Value = "http://172.16.2.61:1616/UI#"
$ws.Range("A3").Value = "http://172.16.2.61:1616/UI#"
$ws.Range("A17").Value = "http://172.16.2.61:1616/UI#"

$ws.Range("A16").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A17").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 3) Row 2: Password + DomainName updated.
# ------------------------------------------------------------------
$ws.Range("C2").Value = "Tetherfi@900"

$ws.Range("E2").Value = "QATetherfi"
$ws.Range("E2").Font.Name = "OCMFont"
$ws.Range("E2").Font.Size = 10
$ws.Range("E2").Font.Color = 4473924

# Row 3: DomainName updated to the same new value + font as E2.
$ws.Range("E3").Value = "QATetherfi"
$ws.Range("E2").Copy()
$ws.Range("E3").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 4) Hyperlinks: A2/A3 no longer carry a live hyperlink (plain text
#    now); C17 (the copied password cell) gets the same mailto
#    hyperlink that C3 originally had.
# ------------------------------------------------------------------
foreach ($hl in @($ws.Hyperlinks)) {
    $addr = $hl.Range.Address()
    if ($addr -eq "`$A`$2" -or $addr -eq "`$A`$3") {
        $hl.Delete()
    }
}

$ws.Hyperlinks.Add($ws.Range("C17"), "mailto:P@ssw0rd@123", "", "", "P@ssw0rd@123")

# ------------------------------------------------------------------
# 5) Data validation: LoginType list now also covers the new D17 cell.
# ------------------------------------------------------------------
$ws.Range("D16").Validation.Delete()
$ws.Range("D16:D17").Validation.Add(3, 1, 1, "(LoginType)")

# ------------------------------------------------------------------
# 6) Selection moves to E3.
# ------------------------------------------------------------------
$ws.Activate()
$ws.Range("E3").Select()
